$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing the existing "Prejuveniles"/"Juveniles"
# rows down by one (old row 11 -> 12, old row 12 -> 13, old row 13 -> 14).
$ws.Rows("11").Insert()

# Fill the newly inserted row 11 with the new "Albatros / caballeros" entry.
$ws.Range("A11").Value = "Torneo FEG"
$ws.Range("B11").Value = "Albatros"
$ws.Range("C11").Value = "caballeros"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "Montoto, Manuel"
$ws.Range("F11").Value = 79

# Append a new row 15 with another "Juveniles / caballeros" entry.
$ws.Range("A15").Value = "Torneo FEG"
$ws.Range("B15").Value = "Juveniles"
$ws.Range("C15").Value = "caballeros"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "Porzio, Tomás"
$ws.Range("F15").Value = 82
